$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.488.20'
$ws.Range('E2').Value = '  -3.04%  '
$ws.Range('D3').Value = '2.374.69'
$ws.Range('E3').Value = '  -4.22%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.36'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.04'
$ws.Range('E6').Value = '  -6.95%  '
$ws.Range('E7').Value = '  -4.57%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('E9').Value = '  -5.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0821'
$ws.Range('E10').Value = '  -4.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '30.97'
$ws.Range('E11').Value = '  -6.79%  '
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').Value = '2.742.13'
$ws.Range('E13').Value = '  -4.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.62'
$ws.Range('E14').Value = '  -4.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.07'
$ws.Range('E15').Value = '  -4.50%  '
$ws.Range('D16').Value = '2.381.76'
$ws.Range('E16').Value = '  -3.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.755'
$ws.Range('E17').Value = '  -4.73%  '
$ws.Range('D18').Value = '40.429.72'
$ws.Range('E18').Value = '  -3.08%  '
$ws.Range('D19').Value = '0.0₃0910'
$ws.Range('E19').Value = '  -4.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.15'
$ws.Range('E20').Value = '  -5.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.88'
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.73'
$ws.Range('E22').Value = '  -5.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.74'
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.63'
$ws.Range('E24').Value = '  -4.48%  '
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.82'
$ws.Range('E26').Value = '  -6.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.55'
$ws.Range('E27').Value = '  -5.04%  '
$ws.Range('E28').Value = '  -2.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.36'
$ws.Range('E29').Value = '  -4.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.74'
$ws.Range('E30').Value = '  -6.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.42'
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.21'
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0728'
$ws.Range('E34').Value = '  -5.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  -6.98%  '
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.81'
$ws.Range('E37').Value = '  -4.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.86'
$ws.Range('E38').Value = '  -8.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0986'
$ws.Range('E39').Value = '  -5.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.73'
$ws.Range('E40').Value = '  -8.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.81'
$ws.Range('E41').Value = '  -5.94%  '
$ws.Range('E42').Value = '  -7.56%  '
$ws.Range('D43').Value = '1.956.94'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0270'
$ws.Range('E44').Value = '  -5.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.61'
$ws.Range('E45').Value = '  -9.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('E46').Value = '  -6.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.30'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').Value = '2.611.19'
$ws.Range('E48').Value = '  -3.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '72.72'
$ws.Range('E49').Value = '  -2.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '93.17'
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.23'
$ws.Range('E51').Value = '  -4.56%  '
